$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.033.99"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "3.270.27"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "'185.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'580.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "3.268.83"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").Value = "'0.131"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.43%  "
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("D12").Value = "'0.412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").Value = "3.839.58"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "'27.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.61%  "
$ws.Range("D16").Value = "68.022.04"
$ws.Range("E16").Value = "  -0.58%  "
$ws.Range("D17").Value = "'0.0000169"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.34%  "
$ws.Range("D18").Value = "3.266.61"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "'5.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("D20").Value = "'13.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "'394.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.73%  "
$ws.Range("D22").Value = "'7.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("D23").Value = "'71.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "'0.513"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  -3.90%  "
$ws.Range("E27").Value = "  -3.30%  "
$ws.Range("D28").Value = "'9.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").Value = "  -2.53%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'5.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.79%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'22.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("E34").Value = "  -4.99%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'163.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  -4.44%  "
$ws.Range("D38").Value = "'1.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("D39").Value = "'26.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "'0.811"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.22%  "
$ws.Range("D41").Value = "'4.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").Value = "'6.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.18%  "
$ws.Range("D43").Value = "'0.0688"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.629.96"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").Value = "'40.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'2.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.88%  "
$ws.Range("D47").Value = "'24.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.08%  "
$ws.Range("D48").Value = "'332.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("D49").Value = "'0.0278"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.70%  "
$ws.Range("D50").Value = "'6.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("E51").Value = "  -1.12%  "
